$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) that look like plain numbers must be forced to
# stay as literal text (matching the sheet's existing text-based layout);
# otherwise Excel would silently convert them to numeric values and lose
# their original formatting/precision.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D16", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D31", "D32", "D33", "D34", "D36", "D39", "D40", "D42", "D43", "D44", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '43.955.65'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '2.296.51'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").Value = '113.97'
$ws.Range("E5").Value = '  +18.60%  '
$ws.Range("D6").Value = '270.75'
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("D7").Value = '0.627'
$ws.Range("E7").Value = '  +0.44%  '
$ws.Range("E8").Value = '  +0.32%  '
$ws.Range("D9").Value = '0.621'
$ws.Range("E9").Value = '  +2.16%  '
$ws.Range("D10").Value = '48.04'
$ws.Range("E10").Value = '  +5.38%  '
$ws.Range("D11").Value = '0.0949'
$ws.Range("E11").Value = '  +1.40%  '
$ws.Range("D12").Value = '9.06'
$ws.Range("E12").Value = '  +13.72%  '
$ws.Range("E13").Value = '  +0.20%  '
$ws.Range("D14").Value = '15.87'
$ws.Range("E14").Value = '  +1.23%  '
$ws.Range("D15").Value = '2.644.28'
$ws.Range("E15").Value = '  +0.37%  '
$ws.Range("D16").Value = '0.854'
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("D17").Value = '2.310.41'
$ws.Range("E17").Value = '  -4.97%  '
$ws.Range("D18").Value = '43.796.35'
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("D19").Value = '0.0000110'
$ws.Range("E19").Value = '  -0.82%  '
$ws.Range("D20").Value = '6.84'
$ws.Range("E20").Value = '  +10.69%  '
$ws.Range("E21").Value = '  +0.70%  '
$ws.Range("D22").Value = '2.44'
$ws.Range("E22").Value = '  -1.23%  '
$ws.Range("B23").Value = 'PancakeSwap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D23").Value = '3.02'
$ws.Range("E23").Value = '  +11.80%  '
$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").Value = '233.18'
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("D25").Value = '9.67'
$ws.Range("E25").Value = '  +6.37%  '
$ws.Range("D26").Value = '11.72'
$ws.Range("E26").Value = '  +3.69%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").Value = '41.89'
$ws.Range("E28").Value = '  +7.91%  '
$ws.Range("D29").Value = '3.40'
$ws.Range("E29").Value = '  -1.86%  '
$ws.Range("E30").Value = '  +1.83%  '
$ws.Range("D31").Value = '175.43'
$ws.Range("E31").Value = '  +0.33%  '
$ws.Range("D32").Value = '0.0935'
$ws.Range("E32").Value = '  +4.39%  '
$ws.Range("D33").Value = '21.62'
$ws.Range("E33").Value = '  -2.00%  '
$ws.Range("D34").Value = '5.71'
$ws.Range("E34").Value = '  +5.20%  '
$ws.Range("E35").Value = '  +0.99%  '
$ws.Range("D36").Value = '4.66'
$ws.Range("E36").Value = '  +1.68%  '
$ws.Range("E37").Value = '  +3.87%  '
$ws.Range("E38").Value = '  +1.31%  '
$ws.Range("D39").Value = '3.81'
$ws.Range("E39").Value = '  +6.67%  '
$ws.Range("D40").Value = '74.17'
$ws.Range("E40").Value = '  +14.98%  '
$ws.Range("E41").Value = '  +3.37%  '
$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").Value = '13.70'
$ws.Range("E42").Value = '  +11.12%  '
$ws.Range("B43").Value = 'THORChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D43").Value = '6.42'
$ws.Range("E43").Value = '  +23.32%  '
$ws.Range("D44").Value = '2.38'
$ws.Range("E44").Value = '  +3.04%  '
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("E46").Value = '  +4.28%  '
$ws.Range("D47").Value = '8.84'
$ws.Range("E47").Value = '  +1.40%  '
$ws.Range("D48").Value = '102.25'
$ws.Range("E48").Value = '  +4.98%  '
$ws.Range("D49").Value = '0.1000'
$ws.Range("E49").Value = '  -2.11%  '
$ws.Range("D50").Value = '0.471'
$ws.Range("E50").Value = '  +9.41%  '
$ws.Range("E51").Value = '  +3.23%  '
